# Fruta / hortaliza, semanal
# Updates the "Fecha", "Calidad", "Volumen", "Precio minimo", "Precio maximo",
# "Precio promedio ponderado" and "Precio $/Kg" columns (D, L, M, N, O, P, S)
# for each data row (2-31) of the weekly Guayaba price sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @{ Row=2; D=44386; L="Primera"; M=160; N=700; O=750; P=725; S=725 },
    @{ Row=3; D=44386; L="Segunda"; M=200; N=600; O=650; P=625; S=625 },
    @{ Row=4; D=44414; L="Primera"; M=160; N=1300; O=1400; P=1350; S=1350 },
    @{ Row=5; D=44351; L="Primera"; M=100; N=700; O=800; P=750; S=750 },
    @{ Row=6; D=44351; L="Segunda"; M=100; N=600; O=700; P=650; S=650 },
    @{ Row=7; D=44260; L="Primera"; M=100; N=1900; O=2000; P=1950; S=1950 },
    @{ Row=8; D=44348; L="Primera"; M=120; N=1000; O=1100; P=1050; S=1050 },
    @{ Row=9; D=44326; L="Primera"; M=160; N=600; O=700; P=650; S=650 },
    @{ Row=10; D=44358; L="Primera"; M=200; N=700; O=800; P=750; S=750 },
    @{ Row=11; D=44358; L="Segunda"; M=200; N=600; O=650; P=625; S=625 },
    @{ Row=12; D=44425; L="Primera"; M=140; N=1200; O=1300; P=1250; S=1250 },
    @{ Row=13; D=44403; L="Primera"; M=100; N=1200; O=1300; P=1250; S=1250 },
    @{ Row=14; D=44403; L="Segunda"; M=120; N=950; O=1000; P=975; S=975 },
    @{ Row=15; D=44344; L="Primera"; M=140; N=1000; O=1200; P=1100; S=1100 },
    @{ Row=16; D=44344; L="Segunda"; M=120; N=800; O=850; P=825; S=825 },
    @{ Row=17; D=44309; L="Primera"; M=160; N=1400; O=1500; P=1450; S=1450 },
    @{ Row=18; D=44417; L="Primera"; M=200; N=1300; O=1400; P=1350; S=1350 },
    @{ Row=19; D=44350; L="Primera"; M=140; N=750; O=800; P=775; S=775 },
    @{ Row=20; D=44372; L="Primera"; M=900; N=750; O=800; P=772; S=772 },
    @{ Row=21; D=44372; L="Segunda"; M=900; N=600; O=650; P=628; S=628 },
    @{ Row=22; D=44498; L="Segunda"; M=100; N=1200; O=1300; P=1250; S=1250 },
    @{ Row=23; D=44407; L="Primera"; M=200; N=600; O=650; P=625; S=625 },
    @{ Row=24; D=44316; L="Primera"; M=140; N=1100; O=1200; P=1150; S=1150 },
    @{ Row=25; D=44389; L="Primera"; M=140; N=750; O=800; P=775; S=775 },
    @{ Row=26; D=44389; L="Segunda"; M=120; N=600; O=700; P=650; S=650 },
    @{ Row=27; D=44330; L="Primera"; M=200; N=1200; O=1300; P=1250; S=1250 },
    @{ Row=28; D=44330; L="Segunda"; M=100; N=1000; O=1100; P=1050; S=1050 },
    @{ Row=29; D=44379; L="Primera"; M=150; N=700; O=800; P=747; S=747 },
    @{ Row=30; D=44379; L="Segunda"; M=140; N=500; O=600; P=543; S=543 },
    @{ Row=31; D=44473; L="Primera"; M=160; N=1500; O=1600; P=1550; S=1550 }

)

foreach ($r in $rows) {
    $ws.Range("D" + $r.Row).Value = $r.D
    $ws.Range("L" + $r.Row).Value = $r.L
    $ws.Range("M" + $r.Row).Value = $r.M
    $ws.Range("N" + $r.Row).Value = $r.N
    $ws.Range("O" + $r.Row).Value = $r.O
    $ws.Range("P" + $r.Row).Value = $r.P
    $ws.Range("S" + $r.Row).Value = $r.S
}
